$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename header text in B1 from "Hourly rate" to "$/hour"
$ws.Range("B1").Value = '$/hour'

# Change number format of B2:B4 from custom currency to plain 2-decimal format
$ws.Range("B2:B4").NumberFormat = "0.00"

# Re-enter the ROUNDUP formulas for C2:E4 as one range-fill operation so
# Excel groups them into a shared formula (matches the diff's t="shared" si="0")
$ws.Range("C2:E4").Formula = '=ROUNDUP($I2*C$6/10,0)*10'

# Update selection to B1 (matches the added <selection activeCell="B1" sqref="B1"/>)
$ws.Range("B1").Select()
